$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet by copying "2021-Q4" (same column
#    layout / header row / styles), placed right after "2021-Q4" and
#    before "总计".
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("2021-Q4")
$ws4Index = $ws4.Index
$ws4.Copy($null, $ws4)
$newSheet = $wb.Worksheets.Item($ws4Index + 1)
$newSheet.Name = "2022-Q1"

# Overwrite the data row with the 2022-Q1 figures. Use a leading
# apostrophe so numeric-looking strings ("002863", "0.44", ...) are
# kept as text (matching the source data) instead of being coerced to
# numbers, while leaving the cell's style/number-format untouched.
$newSheet.Range("B2").Value = "'002863"
$newSheet.Range("C2").Value = "金信深圳成长灵活配置混合"
$newSheet.Range("D2").Value = "'0.44"
$newSheet.Range("E2").Value = "'94.54"
$newSheet.Range("F2").Value = "'6.22"
$newSheet.Range("G2").Value = "'0.0274"
$newSheet.Range("H2").Value = 5

# ---------------------------------------------------------------------
# 2) Prepend a new row to "总计" for the 2022-Q1 summary figures,
#    shifting the existing rows down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Copy A4's format down into the brand-new A5 cell so it picks up the
# same style (bordered/bold index column) as the rest of column A.
$total.Range("A4").Copy($total.Range("A5"))

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.03

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.03

# ---------------------------------------------------------------------
# 3) Restore the originally-active sheet/selection (Copy() above makes
#    the freshly inserted sheet the active one).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
